$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 4) mirroring the structure of existing rows,
# representing a new 20-minute trade entry where the trader could not
# locate close price data from Yahoo (so fewer populated fields than before).
$ws.Range("B4").Value = -12
$ws.Range("C4").Value = 51
$ws.Range("D4").Value = 47
$ws.Range("E4").Value = 18
$ws.Range("F4").Value = 81
$ws.Range("G4").Value = 19774
$ws.Range("H4").Value = 15965
$ws.Range("I4").Value = 951
$ws.Range("J4").Value = 126
$ws.Range("K4").Value = 115
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 9
$ws.Range("N4").Value = "Named"

$ws.Range("A4").Value = 42607.890335648146
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
